$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new data row above row 97 (shifts rows 97-109 down to 98-110),
#    carrying the merged-cell layout / styling down automatically for the
#    pre-existing rows. Excel's plain Insert() drops the thin border that the
#    data rows use, so instead we insert a blank row and then paste the
#    formatting (not the values) from the row that is now directly below the
#    new blank row (the old row 97, now row 98) back onto the blank row.
# ---------------------------------------------------------------------------
$ws.Rows.Item(97).Insert(-4121)

$ws.Range("A98:Q98").Copy()
$ws.Range("A97:Q97").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match row height / merges that the other data rows carry.
$ws.Rows.Item(97).RowHeight = 25.5
$ws.Range("A97:B97").Merge()
$ws.Range("C97:G97").Merge()
$ws.Range("H97:K97").Merge()
$ws.Range("L97:M97").Merge()
$ws.Range("N97:O97").Merge()

# ---------------------------------------------------------------------------
# 2. Fill in the new product row's values. Columns N/P hold numeric-looking
#    text ("450.00", "22.5000") that must stay text (matching the rest of the
#    sheet, which stores every figure here as a shared string), so the number
#    format is forced to Text for the write and then restored so the
#    underlying style index is unaffected.
# ---------------------------------------------------------------------------
$ws.Range("A97").Value = 91
$ws.Range("C97").Value = "حفاضات كبار سن جير ميني 36ق"
$ws.Range("H97").Value = "0:8"

$fmt = $ws.Range("L97").NumberFormat
$ws.Range("L97").NumberFormat = "@"
$ws.Range("L97").Value = "0"
$ws.Range("L97").NumberFormat = $fmt

$fmt = $ws.Range("N97").NumberFormat
$ws.Range("N97").NumberFormat = "@"
$ws.Range("N97").Value = "450.00"
$ws.Range("N97").NumberFormat = $fmt

$fmt = $ws.Range("P97").NumberFormat
$ws.Range("P97").NumberFormat = "@"
$ws.Range("P97").Value = "22.5000"
$ws.Range("P97").NumberFormat = $fmt

$ws.Range("Q97").Value = "0:2"

# ---------------------------------------------------------------------------
# 3. Renumber the "م" (#) column for the rows that followed (they kept their
#    original numbers after the shift; bump each by one so the sequence stays
#    consecutive, 91..102).
# ---------------------------------------------------------------------------
for ($r = 98; $r -le 108; $r++) {
    $cur = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 1).Value = $cur + 1
}

# ---------------------------------------------------------------------------
# 4. Update the running total (now on row 109) to include the new line.
# ---------------------------------------------------------------------------
$ws.Range("P109").Value = 6132.4399999999996

# ---------------------------------------------------------------------------
# 5. Update the generated-at timestamp in the footer (now row 110).
# ---------------------------------------------------------------------------
$ws.Range("A110").Value = "Wednesday, 8 October, 2025 8:32 PM"
